# "hoàn thiện import export kpi item"
# Finalize the KPI item import/export template:
#  - Bold the "Kỳ" label (C2) to match the other bold labels (e.g. "Năm" in E2)
#  - Rename the "Doanh số" header to "Doanh thu"
#  - Simplify the "{Mã NV- Tên nhân viên}" placeholder row into two cells:
#       A5 = "{Mã NV}", B5 = "{Tên nhân viên}" (B5 picks up A5's formatting)
#  - Split the trailing "END - Vui lòng insert ..." text into
#       A10 = "END", B10 = "Vui lòng insert các mã sản phẩm vào theo từng nhân viên"
#       (B10 picks up A10's formatting)

$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KPI nhân viên")

# Bold the "Kỳ" label the same way the "Năm" label (E2) already is.
$ws.Range("E2").Copy()
$ws.Range("C2").PasteSpecial($xlPasteFormats)

# Header row: "Doanh số" -> "Doanh thu"
$ws.Range("E4").Value = "Doanh thu"

# Placeholder row: split combined placeholder into separate NV code / name cells.
# B5 takes on A5's highlighted style before the text is written.
$ws.Range("A5").Copy()
$ws.Range("B5").PasteSpecial($xlPasteFormats)
$ws.Range("A5").Value = "{Mã NV}"
$ws.Range("B5").Value = "{Tên nhân viên}"

# Footer row: split "END - Vui lòng insert..." into END marker + note.
# Write the note (B10) before the marker (A10) so new shared strings land
# in the same order the original authoring tool produced them.
$ws.Range("A10").Copy()
$ws.Range("B10").PasteSpecial($xlPasteFormats)
$ws.Range("B10").Value = "Vui lòng insert các mã sản phẩm vào theo từng nhân viên"
$ws.Range("A10").Value = "END"
